$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit Text
# number format first, otherwise Excel's COM layer auto-converts the
# assigned string into a numeric value (e.g. "0.660" -> 0.66), which
# would lose the exact formatting the source data relies on.
$textCells = @("D5", "D6", "D8", "D9", "D11", "D13", "D15", "D17", "D19", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.380.50'
$ws.Range('E2').Value = '  +3.60%  '
$ws.Range('D3').Value = '2.097.53'
$ws.Range('E3').Value = '  +4.88%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = '250.49'
$ws.Range('E5').Value = '  +2.48%  '
$ws.Range('D6').Value = '0.660'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '51.32'
$ws.Range('E8').Value = '  +16.67%  '
$ws.Range('D9').Value = '61.21'
$ws.Range('E9').Value = '  +9.68%  '
$ws.Range('E10').Value = '  +4.95%  '
$ws.Range('D11').Value = '0.0743'
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('E12').Value = '  +6.88%  '
$ws.Range('D13').Value = '15.27'
$ws.Range('E13').Value = '  +8.15%  '
$ws.Range('D14').Value = '2.401.75'
$ws.Range('E14').Value = '  +4.62%  '
$ws.Range('D15').Value = '0.828'
$ws.Range('E15').Value = '  +5.39%  '
$ws.Range('D16').Value = '2.103.68'
$ws.Range('E16').Value = '  +5.68%  '
$ws.Range('D17').Value = '5.09'
$ws.Range('E17').Value = '  +5.64%  '
$ws.Range('D18').Value = '37.224.18'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('D19').Value = '71.94'
$ws.Range('E19').Value = '  +2.61%  '
$ws.Range('E20').Value = '  +3.19%  '
$ws.Range('D21').Value = '13.58'
$ws.Range('E21').Value = '  +6.40%  '
$ws.Range('D22').Value = '239.86'
$ws.Range('E22').Value = '  +4.02%  '
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').Value = '  +5.76%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').Value = '170.02'
$ws.Range('E26').Value = '  +6.23%  '
$ws.Range('D27').Value = '9.14'
$ws.Range('E27').Value = '  +9.49%  '
$ws.Range('D28').Value = '20.57'
$ws.Range('E28').Value = '  +6.84%  '
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '0.122'
$ws.Range('E30').Value = '  +1.23%  '
$ws.Range('E31').Value = '  +28.61%  '
$ws.Range('D32').Value = '4.48'
$ws.Range('E32').Value = '  +4.13%  '
$ws.Range('E33').Value = '  +6.07%  '
$ws.Range('D34').Value = '0.0920'
$ws.Range('E34').Value = '  +11.34%  '
$ws.Range('B35').Value = 'Gas'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D35').Value = '19.60'
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.33'
$ws.Range('E37').Value = '  +11.51%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = '1.83'
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '4.08'
$ws.Range('E39').Value = '  +2.63%  '
$ws.Range('D40').Value = '1.31'
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').Value = '17.84'
$ws.Range('E41').Value = '  +12.41%  '
$ws.Range('D42').Value = '0.0223'
$ws.Range('E42').Value = '  +5.32%  '
$ws.Range('E43').Value = '  +10.38%  '
$ws.Range('D44').Value = '98.69'
$ws.Range('E44').Value = '  +3.18%  '
$ws.Range('D45').Value = '0.0902'
$ws.Range('E45').Value = '  +12.11%  '
$ws.Range('E46').Value = '  +2.26%  '
$ws.Range('D47').Value = '3.02'
$ws.Range('E47').Value = '  +9.82%  '
$ws.Range('D48').Value = '1.319.08'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').Value = '6.97'
$ws.Range('E49').Value = '  +16.11%  '
$ws.Range('D50').Value = '2.279.96'
$ws.Range('E50').Value = '  +4.17%  '
$ws.Range('D51').Value = '2.28'
$ws.Range('E51').Value = '  +4.82%  '
